$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated data values (refreshed figures from RR GitHub updates) ---
$ws.Range("D5").Value = 70
$ws.Range("E5").Value = 43
$ws.Range("H5").Value = 0.61946902654867253
$ws.Range("I5").Value = 0.38053097345132741
$ws.Range("K5").Value = 5898
$ws.Range("L5").Value = 1827
$ws.Range("O5").Value = 0.76349514563106791
$ws.Range("P5").Value = 0.23650485436893204
$ws.Range("D9").Value = 164
$ws.Range("E9").Value = 82
$ws.Range("H9").Value = 0.66666666666666663
$ws.Range("I9").Value = 0.33333333333333331
$ws.Range("K9").Value = 14420
$ws.Range("L9").Value = 3209
$ws.Range("O9").Value = 0.81797038969879177
$ws.Range("P9").Value = 0.18202961030120823
$ws.Range("D14").Value = 871
$ws.Range("E14").Value = 236
$ws.Range("H14").Value = 0.78681120144534777
$ws.Range("I14").Value = 0.21318879855465223
$ws.Range("K14").Value = 81723
$ws.Range("L14").Value = 10636
$ws.Range("O14").Value = 0.88484067605755801
$ws.Range("P14").Value = 0.11515932394244199
$ws.Range("D15").Value = 654
$ws.Range("E15").Value = 319
$ws.Range("H15").Value = 0.67214799588900309
$ws.Range("I15").Value = 0.32785200411099691
$ws.Range("K15").Value = 58526
$ws.Range("L15").Value = 15850
$ws.Range("O15").Value = 0.78689362159836507
$ws.Range("P15").Value = 0.21310637840163493
$ws.Range("D16").Value = 876
$ws.Range("E16").Value = 208
$ws.Range("H16").Value = 0.80811808118081185
$ws.Range("I16").Value = 0.1918819188191882
$ws.Range("K16").Value = 73242
$ws.Range("L16").Value = 9238
$ws.Range("O16").Value = 0.88799709020368578
$ws.Range("P16").Value = 0.11200290979631426
$ws.Range("D17").Value = 973
$ws.Range("E17").Value = 252
$ws.Range("H17").Value = 0.79428571428571426
$ws.Range("I17").Value = 0.20571428571428571
$ws.Range("K17").Value = 90007
$ws.Range("L17").Value = 12756
$ws.Range("O17").Value = 0.8758697196461761
$ws.Range("P17").Value = 0.12413028035382384
$ws.Range("D18").Value = 1675
$ws.Range("E18").Value = 876
$ws.Range("H18").Value = 0.65660525284202276
$ws.Range("I18").Value = 0.34339474715797724
$ws.Range("K18").Value = 150882
$ws.Range("L18").Value = 43339
$ws.Range("O18").Value = 0.77685729143604454
$ws.Range("P18").Value = 0.22314270856395549
$ws.Range("D19").Value = 597
$ws.Range("E19").Value = 359
$ws.Range("H19").Value = 0.62447698744769875
$ws.Range("I19").Value = 0.37552301255230125
$ws.Range("K19").Value = 53645
$ws.Range("L19").Value = 17346
$ws.Range("O19").Value = 0.75565916806355737
$ws.Range("P19").Value = 0.24434083193644265
$ws.Range("D23").Value = 1071
$ws.Range("E23").Value = 649
$ws.Range("H23").Value = 0.62267441860465111
$ws.Range("I23").Value = 0.37732558139534883
$ws.Range("K23").Value = 97342
$ws.Range("L23").Value = 33717
$ws.Range("O23").Value = 0.74273418841895633
$ws.Range("P23").Value = 0.25726581158104367
$ws.Range("D24").Value = 399
$ws.Range("E24").Value = 191
$ws.Range("H24").Value = 0.67627118644067796
$ws.Range("I24").Value = 0.32372881355932204
$ws.Range("K24").Value = 34621
$ws.Range("L24").Value = 9692
$ws.Range("O24").Value = 0.78128314490104489
$ws.Range("P24").Value = 0.21871685509895517
$ws.Range("D29").Value = 2414
$ws.Range("E29").Value = 696
$ws.Range("H29").Value = 0.77620578778135052
$ws.Range("I29").Value = 0.22379421221864951
$ws.Range("K29").Value = 229388
$ws.Range("L29").Value = 39041
$ws.Range("O29").Value = 0.85455744349530038
$ws.Range("P29").Value = 0.14544255650469956
$ws.Range("D32").Value = 353
$ws.Range("E32").Value = 174
$ws.Range("H32").Value = 0.66982922201138517
$ws.Range("I32").Value = 0.33017077798861483
$ws.Range("K32").Value = 35996
$ws.Range("L32").Value = 8831
$ws.Range("O32").Value = 0.80299819305329379
$ws.Range("P32").Value = 0.19700180694670621
$ws.Range("D33").Value = 649
$ws.Range("E33").Value = 422
$ws.Range("H33").Value = 0.6059757236227824
$ws.Range("I33").Value = 0.39402427637721754
$ws.Range("K33").Value = 55647
$ws.Range("L33").Value = 20742
$ws.Range("O33").Value = 0.72846875859089655
$ws.Range("P33").Value = 0.2715312414091034
$ws.Range("D36").Value = 658
$ws.Range("E36").Value = 163
$ws.Range("H36").Value = 0.80146163215590738
$ws.Range("I36").Value = 0.19853836784409257
$ws.Range("K36").Value = 58524
$ws.Range("L36").Value = 6311
$ws.Range("O36").Value = 0.90266059998457626
$ws.Range("P36").Value = 0.097339400015423766
$ws.Range("D40").Value = 510
$ws.Range("E40").Value = 160
$ws.Range("H40").Value = 0.76119402985074625
$ws.Range("I40").Value = 0.23880597014925373
$ws.Range("K40").Value = 43786
$ws.Range("L40").Value = 5279
$ws.Range("O40").Value = 0.89240803016406811
$ws.Range("P40").Value = 0.10759196983593193
$ws.Range("D41").Value = 1155
$ws.Range("E41").Value = 802
$ws.Range("H41").Value = 0.59018906489524781
$ws.Range("I41").Value = 0.40981093510475219
$ws.Range("K41").Value = 106061
$ws.Range("L41").Value = 40558
$ws.Range("O41").Value = 0.72337827975910352
$ws.Range("P41").Value = 0.27662172024089648
$ws.Range("D44").Value = 1247
$ws.Range("E44").Value = 1269
$ws.Range("H44").Value = 0.49562798092209859
$ws.Range("I44").Value = 0.50437201907790141
$ws.Range("K44").Value = 115038
$ws.Range("L44").Value = 70592
$ws.Range("O44").Value = 0.61971664062920861
$ws.Range("P44").Value = 0.38028335937079139
$ws.Range("D45").Value = 1294
$ws.Range("E45").Value = 805
$ws.Range("H45").Value = 0.61648404001905666
$ws.Range("I45").Value = 0.38351595998094329
$ws.Range("K45").Value = 116968
$ws.Range("L45").Value = 44406
$ws.Range("O45").Value = 0.72482556049921298
$ws.Range("P45").Value = 0.27517443950078702
$ws.Range("D46").Value = 531
$ws.Range("E46").Value = 155
$ws.Range("H46").Value = 0.77405247813411082
$ws.Range("I46").Value = 0.22594752186588921
$ws.Range("K46").Value = 47494
$ws.Range("L46").Value = 6488
$ws.Range("O46").Value = 0.87981178911489011
$ws.Range("P46").Value = 0.12018821088510985
$ws.Range("D47").Value = 588
$ws.Range("E47").Value = 235
$ws.Range("H47").Value = 0.71445929526123941
$ws.Range("I47").Value = 0.28554070473876064
$ws.Range("K47").Value = 50958
$ws.Range("L47").Value = 11712
$ws.Range("O47").Value = 0.81311632359980857
$ws.Range("P47").Value = 0.18688367640019149
$ws.Range("D49").Value = 549
$ws.Range("E49").Value = 128
$ws.Range("H49").Value = 0.81093057607090102
$ws.Range("I49").Value = 0.18906942392909898
$ws.Range("K49").Value = 49105
$ws.Range("L49").Value = 5087
$ws.Range("O49").Value = 0.90613005609684083
$ws.Range("P49").Value = 0.093869943903159142
$ws.Range("D54").Value = 33208
$ws.Range("E54").Value = 13905
$ws.Range("H54").Value = 0.70485853161547773
$ws.Range("I54").Value = 0.29514146838452232
$ws.Range("K54").Value = 3049248
$ws.Range("L54").Value = 691466
$ws.Range("O54").Value = 0.8151513320718986
$ws.Range("P54").Value = 0.18484866792810142
$ws.Range("D57").Value = 173
$ws.Range("E57").Value = 69
$ws.Range("H57").Value = 0.71487603305785119
$ws.Range("I57").Value = 0.28512396694214875
$ws.Range("K57").Value = 15702
$ws.Range("L57").Value = 3284
$ws.Range("O57").Value = 0.82703044348467292
$ws.Range("P57").Value = 0.17296955651532708
$ws.Range("D59").Value = 115
$ws.Range("E59").Value = 48
$ws.Range("H59").Value = 0.70552147239263807
$ws.Range("I59").Value = 0.29447852760736198
$ws.Range("K59").Value = 10452
$ws.Range("L59").Value = 3258
$ws.Range("O59").Value = 0.76236323851203502
$ws.Range("P59").Value = 0.23763676148796498
$ws.Range("D61").Value = 728
$ws.Range("E61").Value = 215
$ws.Range("H61").Value = 0.77200424178154825
$ws.Range("I61").Value = 0.22799575821845175
$ws.Range("K61").Value = 64601
$ws.Range("L61").Value = 10958
$ws.Range("O61").Value = 0.85497425852644948
$ws.Range("P61").Value = 0.14502574147355046
$ws.Range("D89").Value = 112
$ws.Range("E89").Value = 43
$ws.Range("H89").Value = 0.72258064516129028
$ws.Range("I89").Value = 0.27741935483870966
$ws.Range("K89").Value = 10390
$ws.Range("L89").Value = 2086
$ws.Range("O89").Value = 0.83279897403013792
$ws.Range("P89").Value = 0.16720102596986214
$ws.Range("D90").Value = 112
$ws.Range("E90").Value = 43
$ws.Range("H90").Value = 0.72258064516129028
$ws.Range("I90").Value = 0.27741935483870966
$ws.Range("K90").Value = 10390
$ws.Range("L90").Value = 2086
$ws.Range("O90").Value = 0.83279897403013792
$ws.Range("P90").Value = 0.16720102596986214
$ws.Range("D93").Value = 711
$ws.Range("E93").Value = 192
$ws.Range("H93").Value = 0.78737541528239208
$ws.Range("I93").Value = 0.21262458471760798
$ws.Range("K93").Value = 66838
$ws.Range("L93").Value = 8643
$ws.Range("O93").Value = 0.88549436281978244
$ws.Range("P93").Value = 0.11450563718021754
$ws.Range("D105").Value = 114
$ws.Range("E105").Value = 31
$ws.Range("H105").Value = 0.78620689655172415
$ws.Range("I105").Value = 0.21379310344827587
$ws.Range("K105").Value = 10857
$ws.Range("L105").Value = 1748
$ws.Range("O105").Value = 0.86132487108290356
$ws.Range("P105").Value = 0.13867512891709638
$ws.Range("D109").Value = 462
$ws.Range("E109").Value = 67
$ws.Range("H109").Value = 0.87334593572778829
$ws.Range("I109").Value = 0.12665406427221171
$ws.Range("K109").Value = 44125
$ws.Range("L109").Value = 2376
$ws.Range("O109").Value = 0.94890432463817986
$ws.Range("P109").Value = 0.051095675361820178
$ws.Range("D114").Value = 734
$ws.Range("E114").Value = 227
$ws.Range("H114").Value = 0.76378772112382931
$ws.Range("I114").Value = 0.23621227887617066
$ws.Range("K114").Value = 64643
$ws.Range("L114").Value = 10831
$ws.Range("O114").Value = 0.85649362694437825
$ws.Range("P114").Value = 0.14350637305562181
$ws.Range("D116").Value = 380
$ws.Range("E116").Value = 116
$ws.Range("H116").Value = 0.7661290322580645
$ws.Range("I116").Value = 0.23387096774193547
$ws.Range("K116").Value = 32105
$ws.Range("L116").Value = 4088
$ws.Range("O116").Value = 0.88704998204072616
$ws.Range("P116").Value = 0.1129500179592739
$ws.Range("D124").Value = 8509
$ws.Range("E124").Value = 1925
$ws.Range("H124").Value = 0.81550699635806023
$ws.Range("I124").Value = 0.18449300364193982
$ws.Range("K124").Value = 781878
$ws.Range("L124").Value = 87744
$ws.Range("O124").Value = 0.89910098870543753
$ws.Range("P124").Value = 0.10089901129456247
$ws.Range("D134").Value = 320
$ws.Range("E134").Value = 111
$ws.Range("H134").Value = 0.74245939675174011
$ws.Range("I134").Value = 0.25754060324825984
$ws.Range("K134").Value = 29353
$ws.Range("L134").Value = 4558
$ws.Range("O134").Value = 0.86558933679337091
$ws.Range("P134").Value = 0.13441066320662912
$ws.Range("D138").Value = 3746
$ws.Range("E138").Value = 1227
$ws.Range("H138").Value = 0.75326764528453649
$ws.Range("I138").Value = 0.24673235471546351
$ws.Range("K138").Value = 344947
$ws.Range("L138").Value = 55393
$ws.Range("O138").Value = 0.86163511015636707
$ws.Range("P138").Value = 0.1383648898436329
$ws.Range("D139").Value = 47314
$ws.Range("E139").Value = 18699
$ws.Range("H139").Value = 0.71673761228848865
$ws.Range("I139").Value = 0.28326238771151135
$ws.Range("K139").Value = 4342136
$ws.Range("L139").Value = 933801
$ws.Range("O139").Value = 0.82300755297115946
$ws.Range("P139").Value = 0.17699244702884057

# --- Column width tweaks ---
$ws.Columns.Item(4).ColumnWidth = 7.666666666666667
$ws.Columns.Item(12).ColumnWidth = 6.166666666666667

# --- Update the active cell selection on the sheet ---
$ws.Range("O12").Select()
